$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-01 Monday" "2024-04-02 Tuesday"

Replace-Text "61×98=5978" "21×18=378"
Replace-Text "40×63=2520" "27×96=2592"
Replace-Text "55×67=3685" "53×14=742"
Replace-Text "16×60=960" "46×29=1334"
Replace-Text "22×53=1166" "73×81=5913"

Replace-Text "76×82=6232" "30×51=1530"
Replace-Text "16×37=592" "40×12=480"
Replace-Text "68×96=6528" "84×69=5796"
Replace-Text "94×64=6016" "37×84=3108"
Replace-Text "73×11=803" "41×49=2009"

Replace-Text "86×29=2494" "94×59=5546"
Replace-Text "24×30=720" "70×93=6510"
Replace-Text "42×60=2520" "47×98=4606"
Replace-Text "39×43=1677" "31×55=1705"
Replace-Text "68×42=2856" "26×88=2288"

Replace-Text "35×54=1890" "43×92=3956"
Replace-Text "50×23=1150" "34×21=714"
Replace-Text "51×11=561" "39×26=1014"
Replace-Text "59×75=4425" "60×33=1980"
Replace-Text "12×47=564" "29×44=1276"

Replace-Text "56×35=1960" "59×62=3658"
Replace-Text "98×11=1078" "34×12=408"
Replace-Text "44×80=3520" "39×56=2184"
Replace-Text "96×18=1728" "28×84=2352"
Replace-Text "55×29=1595" "95×96=9120"
